$d = $word.ActiveDocument

$pairs = @(
    @("146×9=1314", "124×6=744"),
    @("295×7=2065", "996×3=2988"),
    @("351×6=2106", "321×3=963"),
    @("544×6=3264", "132×7=924"),
    @("160×8=1280", "571×8=4568"),
    @("154×2=308", "357×9=3213"),
    @("441×6=2646", "386×4=1544"),
    @("443×6=2658", "759×6=4554"),
    @("532×4=2128", "624×6=3744"),
    @("540×8=4320", "744×2=1488"),
    @("773×6=4638", "409×3=1227"),
    @("585×6=3510", "918×6=5508"),
    @("139×3=417", "675×9=6075"),
    @("761×8=6088", "356×7=2492"),
    @("899×7=6293", "739×6=4434"),
    @("163×5=815", "631×4=2524"),
    @("144×4=576", "452×9=4068"),
    @("509×6=3054", "402×5=2010"),
    @("108×7=756", "681×8=5448"),
    @("511×9=4599", "417×3=1251"),
    @("472×9=4248", "823×6=4938"),
    @("881×3=2643", "242×3=726"),
    @("174×3=522", "861×9=7749"),
    @("678×6=4068", "596×8=4768")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
